$wb = $excel.ActiveWorkbook

# ALC!row33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 404
$ws.Range("I33").Value = 404
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 404
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -175
$ws.Range("N33").ClearContents()

# ALC!row97
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 200000460
$ws.Range("J97").Value = 200000460
$ws.Range("L97").Value = 600001380
$ws.Range("N97").Value = -600002372

# ALC!row111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1779
$ws.Range("I111").Value = 1315.8572
$ws.Range("K111").Value = 3947.5716
$ws.Range("M111").Value = -880.5715999999998

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2242.5
$ws.Range("I137").Value = 2088
$ws.Range("K137").Value = 6264
$ws.Range("M137").Value = -3714

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 992.7
$ws.Range("I74").Value = 992.7
$ws.Range("K74").Value = 992.7
$ws.Range("M74").Value = -118.7

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 992.7
$ws.Range("I77").Value = 992.7
$ws.Range("K77").Value = 4963.5
$ws.Range("M77").Value = -595.5

# ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2813.125
$ws.Range("I102").Value = 1417.6666
$ws.Range("K102").Value = 1417.6666
$ws.Range("M102").Value = 204.3334

# ARM!row110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2278.8
$ws.Range("I110").Value = 2550
$ws.Range("K110").Value = 2550
$ws.Range("M110").Value = -505

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2999
$ws.Range("I122").Value = 2999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8997
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6547
$ws.Range("N122").ClearContents()

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3758.2
$ws.Range("I132").Value = 3758.2
$ws.Range("K132").Value = 11274.6
$ws.Range("M132").Value = -8744.599999999999

# BSM!row20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2207
$ws.Range("I20").Value = 2207
$ws.Range("K20").Value = 2207
$ws.Range("M20").Value = -1960

# BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3903.4546
$ws.Range("I105").Value = 3492.375
$ws.Range("K105").Value = 3492.375
$ws.Range("M105").Value = -1745.375

# BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2574.875
$ws.Range("I107").Value = 2574.875
$ws.Range("K107").Value = 2574.875
$ws.Range("M107").Value = -654.875

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4994.3335
$ws.Range("I134").Value = 4989
$ws.Range("J134").Value = 4997
$ws.Range("K134").Value = 14967
$ws.Range("L134").Value = 14991
$ws.Range("M134").Value = -12432
$ws.Range("N134").Value = -20061

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2661.6667
$ws.Range("I31").Value = 2947.2
$ws.Range("J31").Value = 1234
$ws.Range("K31").Value = 2947.2
$ws.Range("L31").Value = 1234
$ws.Range("M31").Value = -2652.2
$ws.Range("N31").Value = -1824

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2661.6667
$ws.Range("I34").Value = 2947.2
$ws.Range("J34").Value = 1234
$ws.Range("K34").Value = 2947.2
$ws.Range("L34").Value = 1234
$ws.Range("M34").Value = -2745.2
$ws.Range("N34").Value = -1638

# CRP!row86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11620328
$ws.Range("I86").Value = 17428492
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 17428492
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -17427369
$ws.Range("N86").Value = -6246

# CRP!row89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 11620328
$ws.Range("I89").Value = 17428492
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 87142460
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -87136844
$ws.Range("N89").Value = -31232

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5655.857
$ws.Range("I99").Value = 5103.7646
$ws.Range("K99").Value = 5103.7646
$ws.Range("M99").Value = -3605.7646

# CRP!row105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1860
$ws.Range("I105").Value = 900
$ws.Range("K105").Value = 900
$ws.Range("M105").Value = 847

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1271.2727
$ws.Range("I122").Value = 1930
$ws.Range("J122").Value = 1124.8889
$ws.Range("K122").Value = 5790
$ws.Range("L122").Value = 3374.6667
$ws.Range("M122").Value = -3340
$ws.Range("N122").Value = -8274.6667

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5655.857
$ws.Range("I126").Value = 5103.7646
$ws.Range("K126").Value = 15311.2938
$ws.Range("M126").Value = -12841.2938

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7177.8
$ws.Range("I132").Value = 5248.75
$ws.Range("K132").Value = 15746.25
$ws.Range("M132").Value = -13216.25

# CUL!row5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 376.6
$ws.Range("I5").Value = 292
$ws.Range("J5").Value = 433
$ws.Range("K5").Value = 876
$ws.Range("L5").Value = 1299
$ws.Range("M5").Value = -764
$ws.Range("N5").Value = -1523

# CUL!row18
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1712.8572
$ws.Range("I18").Value = 995
$ws.Range("K18").Value = 2985
$ws.Range("M18").Value = -2816

# CUL!row68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 598.5
$ws.Range("I68").Value = 598
$ws.Range("K68").Value = 1794
$ws.Range("M68").Value = -983

# CUL!row71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 598.5
$ws.Range("I71").Value = 598
$ws.Range("K71").Value = 5382
$ws.Range("M71").Value = -1326

# CUL!row135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 376.6
$ws.Range("I135").Value = 292
$ws.Range("J135").Value = 433
$ws.Range("K135").Value = 2628
$ws.Range("L135").Value = 3897
$ws.Range("M135").Value = -93
$ws.Range("N135").Value = -8967

# GSM!row101
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 35999.668
$ws.Range("J101").Value = 35999.668
$ws.Range("L101").Value = 35999.668
$ws.Range("N101").Value = -42489.668

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3509.5881
$ws.Range("I132").Value = 3482.125
$ws.Range("K132").Value = 10446.375
$ws.Range("M132").Value = -7916.375

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4991.8335
$ws.Range("J7").Value = 4987.5
$ws.Range("L7").Value = 4987.5
$ws.Range("N7").Value = -5211.5

# LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2500
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

# LTW!row25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2500
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

# LTW!row40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3099.25
$ws.Range("I40").Value = 3099.25
$ws.Range("K40").Value = 3099.25
$ws.Range("M40").Value = -2963.25

# LTW!row47
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 500
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# LTW!row52
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H52").Value = 500
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

# LTW!row101
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 44982.668
$ws.Range("J101").Value = 44982.668
$ws.Range("L101").Value = 44982.668
$ws.Range("N101").Value = -51472.668

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4991.8335
$ws.Range("J126").Value = 4987.5
$ws.Range("L126").Value = 14962.5
$ws.Range("N126").Value = -19902.5

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3428.0833
$ws.Range("I132").Value = 3376.0908
$ws.Range("K132").Value = 10128.2724
$ws.Range("M132").Value = -7598.2724

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4513.143
$ws.Range("I136").Value = 3635.625
$ws.Range("K136").Value = 10906.875
$ws.Range("M136").Value = -8356.875

# WVR!row2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 20000
$ws.Range("I2").Value = 20000
$ws.Range("K2").Value = 20000
$ws.Range("M2").Value = -19888

# WVR!row87
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# WVR!row90
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# WVR!row103
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 20499.666
$ws.Range("J103").Value = 20499.666
$ws.Range("L103").Value = 20499.666
$ws.Range("N103").Value = -22843.666

# WVR!row107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 465.9
$ws.Range("J107").Value = 294
$ws.Range("L107").Value = 882
$ws.Range("N107").Value = -4722

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2107.875
$ws.Range("I136").Value = 1461.32
$ws.Range("K136").Value = 4383.96
$ws.Range("M136").Value = -1833.96
